$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StringLocalizations_Valencia")

# Insert a new row above the current row 41 ("BASIC_TEXT_DESCRIPTION" row),
# shifting all subsequent rows down by one.
$ws.Rows.Item(41).Insert()

# Populate the new row with the download-prompt strings (game over screen art asset text).
$ws.Cells.Item(41, 1).Value = "BASIC_TEXT_DOWNLOAD"
$ws.Cells.Item(41, 2).Value = "Download the INSPEC2T App Now!"
$ws.Cells.Item(41, 3).Value = "XXXX"
$ws.Cells.Item(41, 4).Value = "XXXX"
$ws.Cells.Item(41, 5).Value = "XXXX"

# Move the selection to reflect where the editor was last working.
$ws.Activate()
$ws.Range("B43").Select()
